$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 222 (Excel semantics: existing
# row 222 and everything below it shifts down by one).
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new weekly record.
$ws.Cells.Item(222, 1).Value = 3
$ws.Cells.Item(222, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(222, 3).Value = "Coquimbo"
$ws.Cells.Item(222, 4).Value = 44516
$ws.Cells.Item(222, 5).Value = 5
$ws.Cells.Item(222, 6).Value = 100114013
$ws.Cells.Item(222, 7).Value = "Zanahoria"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 540
$ws.Cells.Item(222, 11).Value = 6000
$ws.Cells.Item(222, 12).Value = 6500
$ws.Cells.Item(222, 13).Value = 6259
$ws.Cells.Item(222, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(222, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(222, 16).Value = 313
$ws.Cells.Item(222, 17).Value = 20
$ws.Cells.Item(222, 18).Value = "Hortaliza"
